$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new gauge row (Poamoho rain) with Min_val=0, Max_val=7
$ws.Range("A6").Value = "Poamoho rain"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 7

# Update the current selection to match the target state
$ws.Range("C10").Select()
